$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.253221392631531
$ws.Range("B1").Value = 2.701875448226929
$ws.Range("C1").Value = 8.473056793212891
$ws.Range("D1").Value = 2.081536054611206
$ws.Range("E1").Value = 1.137606263160706
